$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 118
$ws.Range("I2").Value = 390
$ws.Range("J2").Value = 1568
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 430
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = 248
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 167
$ws.Range("T2").Value = 262
$ws.Range("V2").Value = 2309
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2297
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = 19
